# Catalog.xlsx: refresh the product list with the new placeholder
# copy and drop the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Описание" (C) text for the rows that keep their existing
# "Название" (B) value.
$ws.Range("C2").Value = "Ipsum dolor sit amet dolor sit amet dolor sit amet dolor sit amet"
$ws.Range("C3").Value = "Ipsum dolor sit amet"

# Rows 4-8 are repurposed to the next items that used to live further
# down the (now removed) tail of the table.
$ws.Range("B4").Value = "Нарты спортивные"
$ws.Range("C4").Value = "Ipsum dolor sit amet"

$ws.Range("B5").Value = "Сумка для снаряжения"
$ws.Range("C5").Value = "Ipsum dolor sit amet dolor sit amet dolor sit amet dolor sit amet dolor sit amet "

$ws.Range("B6").Value = "Дождевик"
$ws.Range("C6").Value = "Ipsum dolor sit amet"

$ws.Range("B7").Value = "Попона"
$ws.Range("C7").Value = "Ipsum dolor sit amet"

$ws.Range("B8").Value = "Тапочки"
$ws.Range("C8").Value = "Ipsum dolor sit amet dolor sit amet dolor sit ametdolor sit ametdolor sit ametdolor sit ametdolor sit ametdolor sit amet"

# The remaining catalog entries (old rows 9-16) are gone entirely.
$ws.Range("A9:C16").ClearContents()

# Column widths were re-fit to the shorter descriptions.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
